$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 34.97976933333334
$ws.Cells.Item(2, 14).Value = 104.939308
$ws.Cells.Item(2, 15).Value = 0.4352965780925344
$ws.Cells.Item(2, 16).Value = 0.4352965780925344
$ws.Cells.Item(2, 17).Value = 16.70680423052445
$ws.Cells.Item(2, 18).Value = 150.36123807472
$ws.Cells.Item(2, 19).Value = 0.4352965780925344
$ws.Cells.Item(2, 20).Value = 0.4352965780925344

# Row 3 updates
$ws.Cells.Item(3, 14).Value = 61.03014900000001
$ws.Cells.Item(3, 15).Value = 0.2531579017099818
$ws.Cells.Item(3, 16).Value = 0.2531579017099818
$ws.Cells.Item(3, 17).Value = 9.71627096590667
$ws.Cells.Item(3, 18).Value = 87.44643869316002
$ws.Cells.Item(3, 19).Value = 0.2531579017099818
$ws.Cells.Item(3, 20).Value = 0.2531579017099818

# Row 4 updates
$ws.Cells.Item(4, 13).Value = 25.035323
$ws.Cells.Item(4, 14).Value = 75.105969
$ws.Cells.Item(4, 15).Value = 0.3115455201974837
$ws.Cells.Item(4, 16).Value = 0.3115455201974837
$ws.Cells.Item(4, 17).Value = 11.95720406910667
$ws.Cells.Item(4, 18).Value = 107.61483662196
$ws.Cells.Item(4, 19).Value = 0.3115455201974837
$ws.Cells.Item(4, 20).Value = 0.3115455201974837

# Remove row 5 entirely (shifts dimension to A1:T4 and drops the now-unused
# "Resolving-Mac" shared string)
$ws.Rows.Item(5).Delete()
